# Update the zuper_id for the "12B234" / Dell machine register row (row 2)
# from 12974 to 5342.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = 5342
